$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1. Rename the three existing category labels
# -----------------------------------------------------------------
$ws.Range("A1").Value = "Total CRS"
$ws.Rows(1).AutoFit()

$ws.Range("A5").Value = "Total SRS"
$ws.Rows(5).AutoFit()

$ws.Range("A9").Value = "Total TC"
$ws.Rows(9).AutoFit()

# -----------------------------------------------------------------
# 2. Add three new blocks (rows 13-24), each a 4-row tall section
#    mirroring the existing A1:F4 / A5:F8 / A9:F12 layout:
#      - A:C merged, holding a header label (filled, wrapped)
#      - D:F merged, holding the (currently empty) value cell
#    Cells are merged BEFORE the formatting is pasted onto them so
#    the merge doesn't fragment the border styling cell-by-cell.
# -----------------------------------------------------------------
$newLabels = @("Total TC`nExecuted", "Total TC`nPassed", "Total TC `nFailed")
$startRows = @(13, 17, 21)

for ($i = 0; $i -lt 3; $i++) {
    $r = $startRows[$i]
    $rEnd = $r + 3

    $ws.Range("A$r`:C$rEnd").Merge()
    $ws.Range("D$r`:F$rEnd").Merge()

    # Clone the formatting of the previous block (A9:F12) onto the new block
    $ws.Range("A9:F12").Copy()
    $ws.Range("A$r`:F$rEnd").PasteSpecial(-4122)

    # Header text + wrap for the label cell
    $ws.Range("A$r").Value = $newLabels[$i]
    $ws.Range("A$r").WrapText = $true
    $ws.Rows($r).AutoFit()
}

# -----------------------------------------------------------------
# 3. Update the view: top-left visible row + selection
# -----------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("D17:F20").Select()
